# Removed 0 for missing ratings.

$wb = $excel.ActiveWorkbook

$wsNotes = $wb.Worksheets.Item("DBNotes")

# Clear the cells that only held a placeholder 0 rating (missing rating)
$zeroCells = @("H14", "I14", "J14", "N14", "G15", "H15", "I15", "J15", "K15", "L15", "N15", "H16", "J16", "K16", "G28", "I28", "J28", "L28", "H31", "G43", "H43", "I43", "M43", "N43", "H46", "G59", "H59", "I59", "M59", "N59", "H61", "G73", "H73", "I73", "L73", "M73", "N73", "I74", "G76", "H76", "G88", "H88", "I88", "L88", "M88", "N88", "I89", "N89", "G91", "H91", "I104", "H106", "G118", "H118", "I118", "M118", "N118", "I119", "H121")

foreach ($addr in $zeroCells) {
    $wsNotes.Range($addr).ClearContents()
}

# Make DBNotes the active sheet/tab, and select columns G:N on it
$wsNotes.Activate()
$wsNotes.Range("G1:N1048576").Select()
